$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Fix price values (24hr change and trend fixed)
$ws.Range("D3").Value = 100
$ws.Range("D7").Value = 40
$ws.Range("D8").Value = 40
$ws.Range("D9").Value = 40

# Move selection to D10 as last active cell
$ws.Range("D10").Select()
